# ============================================================
# Edit: add 2022-Q3 fund-holdings data
#
#  1. Insert a new "2022-Q3" worksheet right after the "总计"
#     (summary) sheet, populated with that quarter's fund table.
#  2. Update the "总计" summary sheet: insert a new row at the top
#     of the data (row 2) holding the 2022-Q3 totals, shifting all
#     the existing quarter rows down by one.
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Insert a blank row at row 2 - this shifts existing rows 2-8 down to
# rows 3-9, carrying their formatting (including column A's style)
# along with them. The engine carries the header row's bold/no-border
# format down onto the new row (matching Excel's own "format same as
# above" default on row-insert), so B2:D2 need their format reset back
# to the plain/default style used by the rest of the data rows.
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# The freshly inserted row 2 has no formatting of its own for column A;
# copy the format for column A from the row directly below (which now
# holds the old row 2's formatting) so A2 matches the rest of the index
# column.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

# Fill in the new 2022-Q3 summary row.
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 25
$summary.Cells.Item(2,4).Value = 1.82

# Column A is a simple 0-based row counter; after inserting the new
# row it needs renumbering down the rest of the table.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(9,1).Value = 7

# ------------------------------------------------------------------
# Step 2: create the new "2022-Q3" worksheet right after "总计".
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (column A has no header, matches the other quarter sheets).
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Column A (row index, 0-based) - always numeric.
for ($i = 0; $i -lt 25; $i++) {
    $q3.Cells.Item($i + 2, 1).Value = $i
}

# Columns B-G hold fund code / name / size / position% / ratio / market
# value - these are stored as TEXT in the source data (even though most
# look numeric), so force a text number format before writing them.
$q3.Range("B2:G26").NumberFormat = "@"

$arrBH = New-Object 'object[,]' 25,7
$arrBH[0,0] = "161219"
$arrBH[0,1] = "国投瑞银新兴产业混合（LOF）"
$arrBH[0,2] = "6.18"
$arrBH[0,3] = "79.94"
$arrBH[0,4] = "4.14"
$arrBH[0,5] = "0.2559"
$arrBH[0,6] = 7
$arrBH[1,0] = "161232"
$arrBH[1,1] = "国投瑞银瑞盛灵活配置混合"
$arrBH[1,2] = "4.18"
$arrBH[1,3] = "94.55"
$arrBH[1,4] = "5.50"
$arrBH[1,5] = "0.2299"
$arrBH[1,6] = 6
$arrBH[2,0] = "000663"
$arrBH[2,1] = "国投瑞银美丽中国灵活配置混合"
$arrBH[2,2] = "3.85"
$arrBH[2,3] = "93.40"
$arrBH[2,4] = "4.82"
$arrBH[2,5] = "0.1856"
$arrBH[2,6] = 6
$arrBH[3,0] = "000968"
$arrBH[3,1] = "广发中证养老产业指数A"
$arrBH[3,2] = "10.84"
$arrBH[3,3] = "93.98"
$arrBH[3,4] = "1.66"
$arrBH[3,5] = "0.1799"
$arrBH[3,6] = 2
$arrBH[4,0] = "001195"
$arrBH[4,1] = "工银农业产业股票"
$arrBH[4,2] = "5.53"
$arrBH[4,3] = "80.67"
$arrBH[4,4] = "2.97"
$arrBH[4,5] = "0.1642"
$arrBH[4,6] = 9
$arrBH[5,0] = "001320"
$arrBH[5,1] = "工银丰盈回报灵活配置混合A"
$arrBH[5,2] = "2.81"
$arrBH[5,3] = "82.02"
$arrBH[5,4] = "4.31"
$arrBH[5,5] = "0.1211"
$arrBH[5,6] = 7
$arrBH[6,0] = "161225"
$arrBH[6,1] = "国投瑞银瑞盈灵活配置混合（LOF）"
$arrBH[6,2] = "2.52"
$arrBH[6,3] = "94.35"
$arrBH[6,4] = "4.79"
$arrBH[6,5] = "0.1207"
$arrBH[6,6] = 6
$arrBH[7,0] = "010994"
$arrBH[7,1] = "博时创新经济混合A"
$arrBH[7,2] = "3.35"
$arrBH[7,3] = "91.74"
$arrBH[7,4] = "3.30"
$arrBH[7,5] = "0.1106"
$arrBH[7,6] = 10
$arrBH[8,0] = "012202"
$arrBH[8,1] = "中加消费优选混合A"
$arrBH[8,2] = "3.52"
$arrBH[8,3] = "73.15"
$arrBH[8,4] = "2.95"
$arrBH[8,5] = "0.1038"
$arrBH[8,6] = 9
$arrBH[9,0] = "013347"
$arrBH[9,1] = "工银丰盈回报灵活配置混合C"
$arrBH[9,2] = "1.22"
$arrBH[9,3] = "82.02"
$arrBH[9,4] = "4.31"
$arrBH[9,5] = "0.0526"
$arrBH[9,6] = 7
$arrBH[10,0] = "000556"
$arrBH[10,1] = "国投瑞银新机遇灵活配置混合A"
$arrBH[10,2] = "4.50"
$arrBH[10,3] = "21.12"
$arrBH[10,4] = "1.04"
$arrBH[10,5] = "0.0468"
$arrBH[10,6] = 8
$arrBH[11,0] = "002358"
$arrBH[11,1] = "国投瑞银瑞祥灵活配置混合A"
$arrBH[11,2] = "4.17"
$arrBH[11,3] = "20.72"
$arrBH[11,4] = "1.11"
$arrBH[11,5] = "0.0463"
$arrBH[11,6] = 7
$arrBH[12,0] = "161233"
$arrBH[12,1] = "国投瑞银瑞泰多策略灵活配置混合（LOF）A"
$arrBH[12,2] = "4.96"
$arrBH[12,3] = "28.33"
$arrBH[12,4] = "0.82"
$arrBH[12,5] = "0.0407"
$arrBH[12,6] = 9
$arrBH[13,0] = "011616"
$arrBH[13,1] = "国投瑞银瑞祥灵活配置混合C"
$arrBH[13,2] = "3.64"
$arrBH[13,3] = "20.72"
$arrBH[13,4] = "1.11"
$arrBH[13,5] = "0.0404"
$arrBH[13,6] = 7
$arrBH[14,0] = "015056"
$arrBH[14,1] = "百嘉百盛混合"
$arrBH[14,2] = "1.21"
$arrBH[14,3] = "62.90"
$arrBH[14,4] = "3.10"
$arrBH[14,5] = "0.0375"
$arrBH[14,6] = 3
$arrBH[15,0] = "012203"
$arrBH[15,1] = "中加消费优选混合C"
$arrBH[15,2] = "0.57"
$arrBH[15,3] = "73.15"
$arrBH[15,4] = "2.95"
$arrBH[15,5] = "0.0168"
$arrBH[15,6] = 9
$arrBH[16,0] = "002982"
$arrBH[16,1] = "广发中证养老产业指数C"
$arrBH[16,2] = "0.92"
$arrBH[16,3] = "93.98"
$arrBH[16,4] = "1.66"
$arrBH[16,5] = "0.0153"
$arrBH[16,6] = 2
$arrBH[17,0] = "000557"
$arrBH[17,1] = "国投瑞银新机遇灵活配置混合C"
$arrBH[17,2] = "1.45"
$arrBH[17,3] = "21.12"
$arrBH[17,4] = "1.04"
$arrBH[17,5] = "0.0151"
$arrBH[17,6] = 8
$arrBH[18,0] = "516560"
$arrBH[18,1] = "华宝养老ETF"
$arrBH[18,2] = "0.74"
$arrBH[18,3] = "98.01"
$arrBH[18,4] = "1.73"
$arrBH[18,5] = "0.0128"
$arrBH[18,6] = 2
$arrBH[19,0] = "010995"
$arrBH[19,1] = "博时创新经济混合C"
$arrBH[19,2] = "0.38"
$arrBH[19,3] = "91.74"
$arrBH[19,4] = "3.30"
$arrBH[19,5] = "0.0125"
$arrBH[19,6] = 10
$arrBH[20,0] = "011618"
$arrBH[20,1] = "国投瑞银瑞泰多策略灵活配置混合（LOF）C"
$arrBH[20,2] = "0.80"
$arrBH[20,3] = "28.33"
$arrBH[20,4] = "0.82"
$arrBH[20,5] = "0.0066"
$arrBH[20,6] = 9
$arrBH[21,0] = "013072"
$arrBH[21,1] = "泰信医疗服务混合A"
$arrBH[21,2] = "0.08"
$arrBH[21,3] = "91.16"
$arrBH[21,4] = "3.38"
$arrBH[21,5] = "0.0027"
$arrBH[21,6] = 8
$arrBH[22,0] = "001657"
$arrBH[22,1] = "长安鑫富领先灵活配置混合"
$arrBH[22,2] = "0.06"
$arrBH[22,3] = "49.60"
$arrBH[22,4] = "3.27"
$arrBH[22,5] = "0.0020"
$arrBH[22,6] = 4
$arrBH[23,0] = "003366"
$arrBH[23,1] = "浙商汇金中证转型成长指数"
$arrBH[23,2] = "0.07"
$arrBH[23,3] = "93.03"
$arrBH[23,4] = "1.27"
$arrBH[23,5] = "0.0009"
$arrBH[23,6] = 4
$arrBH[24,0] = "013073"
$arrBH[24,1] = "泰信医疗服务混合C"
$arrBH[24,2] = "0.01"
$arrBH[24,3] = "91.16"
$arrBH[24,4] = "3.38"
$arrBH[24,5] = "0.0003"
$arrBH[24,6] = 8

$q3.Range("B2:H26").Value = $arrBH

Write-Output "Edit applied successfully."
